$d = $word.ActiveDocument

# 1. Update the FIRST bilibili video link (both the hidden HYPERLINK field
#    code and its displayed text) to the new video id. The document
#    contains the same link twice; only the first occurrence changes.
$oldId = "BV1GFKBzbE7w"
$newId = "BV1LhKfzZEsY"

$targetField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields.Item($i)
    if ($candidate.Data -like "*$oldId*") {
        $targetField = $candidate
        break
    }
}
$targetField.Data = $targetField.Data.Replace($oldId, $newId)

# Locate the paragraph that holds this field's displayed result so the
# text replacement below only touches that one occurrence.
$fieldStart = $targetField.Result.Start
$linkParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($fieldStart -ge $p.Range.Start -and $fieldStart -lt $p.Range.End) {
        $linkParaIndex = $i
        break
    }
}
$linkParagraph = $d.Paragraphs.Item($linkParaIndex)
$linkParagraph.Range.Find.Execute($oldId, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newId, 2)

# 2. Move the hidden "_GoBack" bookmark from its old spot (right before the
#    second hyperlink, near the end of the document) onto the now-empty
#    paragraph right after the first hyperlink. Adding a bookmark with the
#    same reserved name elsewhere automatically removes the previous one,
#    matching how Word keeps only a single "_GoBack" location.
$targetParagraph = $d.Paragraphs.Item($linkParaIndex + 1)
$d.Bookmarks.Add("_GoBack", $targetParagraph.Range)

# 3. Remove the extra duplicate empty paragraph that used to directly
#    follow the now-bookmarked paragraph.
$extraParagraph = $d.Paragraphs.Item($linkParaIndex + 2)
$extraParagraph.Range.Delete()
